$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 10
$ws.Cells.Item(117, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(117, 3).Value = "La Araucanía"
$ws.Cells.Item(117, 4).Value = 44452
$ws.Cells.Item(117, 5).Value = 9
$ws.Cells.Item(117, 6).Value = 100112044
$ws.Cells.Item(117, 7).Value = "Perejil"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 80
$ws.Cells.Item(117, 11).Value = 4000
$ws.Cells.Item(117, 12).Value = 4000
$ws.Cells.Item(117, 13).Value = 4000
$ws.Cells.Item(117, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(117, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(117, 16).Value = 1333
$ws.Cells.Item(117, 17).Value = 3
$ws.Cells.Item(117, 18).Value = "Hortaliza"
